# update title and create gif
# Applies the textual edits to the two body paragraphs and relocates the
# "_GoBack" bookmark from the paragraph after the "Thank you" sentence to
# the paragraph right after the "...if possible" sentence.

$d = $word.ActiveDocument

# --- Paragraph 1 edits (". I will consist of react-strap ... possible.") ---

# "react-strap" -> "reactstrap"
$d.Content.Find.Execute("react-strap", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "reactstrap", 2) | Out-Null

# ". I will" -> ". It will"
$d.Content.Find.Execute(". I will", $false, $false, $false, $false, $false, `
                         $true, 1, $false, ". It will", 2) | Out-Null

# "card components" -> "table components"
$d.Content.Find.Execute("card components", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "table components", 2) | Out-Null

# "API if possible." -> "API if possible, and axios."
$d.Content.Find.Execute("API if possible.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "API if possible, and axios.", 2) | Out-Null

# --- Paragraph 2 edits ("The structure will be ... Thank you") ---

# "type left" -> "top left"
$d.Content.Find.Execute("type left", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "top left", 2) | Out-Null

# "possibly transfer into" -> "possibly transfer this project into"
$d.Content.Find.Execute("possibly transfer into", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "possibly transfer this project into", 2) | Out-Null

# --- Relocate the _GoBack bookmark ---
# It currently sits in the empty paragraph right after the "Thank you"
# paragraph; move it to the empty paragraph right after the
# "...if possible, and axios." paragraph.

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$targetPara = $null
for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*and axios.*") {
        $targetPara = $d.Paragraphs($i + 1)
        break
    }
}

if ($targetPara -ne $null) {
    $gbRange = $targetPara.Range
    $gbRange.Collapse(1)
    $d.Bookmarks.Add("_GoBack", $gbRange)
}
